$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 308.5
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 378
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 378
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -718
$ws.Range("H28").Value = 805.2381
$ws.Range("I28").Value = 979.2353000000001
$ws.Range("J28").Value = 65.75
$ws.Range("K28").Value = 979.2353000000001
$ws.Range("L28").Value = 65.75
$ws.Range("M28").Value = -494.2353000000001
$ws.Range("N28").Value = -1035.75
$ws.Range("H76").Value = 150002450
$ws.Range("I76").Value = 166669000
$ws.Range("K76").Value = 166669000
$ws.Range("M76").Value = -166668685
$ws.Range("H79").Value = 150002450
$ws.Range("I79").Value = 166669000
$ws.Range("K79").Value = 166669000
$ws.Range("M79").Value = -166667908
$ws.Range("H86").Value = 111116030
$ws.Range("I86").Value = 4300.6
$ws.Range("J86").Value = 250005700
$ws.Range("K86").Value = 4300.6
$ws.Range("L86").Value = 250005700
$ws.Range("M86").Value = -3177.6
$ws.Range("N86").Value = -250007946
$ws.Range("H89").Value = 111116030
$ws.Range("I89").Value = 4300.6
$ws.Range("J89").Value = 250005700
$ws.Range("K89").Value = 21503
$ws.Range("L89").Value = 1250028500
$ws.Range("M89").Value = -15887
$ws.Range("N89").Value = -1250039732
$ws.Range("H132").Value = 2780049.5
$ws.Range("I132").Value = 2421.5186
$ws.Range("J132").Value = 11112933
$ws.Range("K132").Value = 7264.5558
$ws.Range("L132").Value = 33338799
$ws.Range("M132").Value = -4734.5558
$ws.Range("N132").Value = -33343859

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2582.85
$ws.Range("I32").Value = 2582.85
$ws.Range("K32").Value = 2582.85
$ws.Range("M32").Value = -2295.85
$ws.Range("H61").Value = 1369.6744
$ws.Range("I61").Value = 1384.9
$ws.Range("J61").Value = 1166.6666
$ws.Range("K61").Value = 1384.9
$ws.Range("L61").Value = 1166.6666
$ws.Range("M61").Value = -1172.9
$ws.Range("N61").Value = -1590.6666
$ws.Range("H74").Value = 1633.3914
$ws.Range("I74").Value = 1673.8
$ws.Range("J74").Value = 1364
$ws.Range("K74").Value = 1673.8
$ws.Range("L74").Value = 1364
$ws.Range("M74").Value = -799.8
$ws.Range("N74").Value = -3112
$ws.Range("H77").Value = 1633.3914
$ws.Range("I77").Value = 1673.8
$ws.Range("J77").Value = 1364
$ws.Range("K77").Value = 8369
$ws.Range("L77").Value = 6820
$ws.Range("M77").Value = -4001
$ws.Range("N77").Value = -15556
$ws.Range("H102").Value = 6458.8237
$ws.Range("I102").Value = 4986.6665
$ws.Range("J102").Value = 17500
$ws.Range("K102").Value = 4986.6665
$ws.Range("L102").Value = 17500
$ws.Range("M102").Value = -3364.6665
$ws.Range("N102").Value = -20744
$ws.Range("H122").Value = 889.5
$ws.Range("I122").Value = 889.5
$ws.Range("K122").Value = 2668.5
$ws.Range("M122").Value = -218.5
$ws.Range("H132").Value = 1752.4615
$ws.Range("I132").Value = 1523.1428
$ws.Range("J132").Value = 2715.6
$ws.Range("K132").Value = 4569.428400000001
$ws.Range("L132").Value = 8146.799999999999
$ws.Range("M132").Value = -2039.428400000001
$ws.Range("N132").Value = -13206.8
$ws.Range("H136").Value = 1369.6744
$ws.Range("I136").Value = 1384.9
$ws.Range("J136").Value = 1166.6666
$ws.Range("K136").Value = 4154.700000000001
$ws.Range("L136").Value = 3499.9998
$ws.Range("M136").Value = -1604.700000000001
$ws.Range("N136").Value = -8599.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 10525.429
$ws.Range("J107").Value = 67400
$ws.Range("L107").Value = 67400
$ws.Range("N107").Value = -71240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43497.457
$ws.Range("I31").Value = 1586.2632
$ws.Range("J31").Value = 202760
$ws.Range("K31").Value = 1586.2632
$ws.Range("L31").Value = 202760
$ws.Range("M31").Value = -1291.2632
$ws.Range("N31").Value = -203350
$ws.Range("H34").Value = 43497.457
$ws.Range("I34").Value = 1586.2632
$ws.Range("J34").Value = 202760
$ws.Range("K34").Value = 1586.2632
$ws.Range("L34").Value = 202760
$ws.Range("M34").Value = -1384.2632
$ws.Range("N34").Value = -203164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3096118.2
$ws.Range("I2").Value = 4525063
$ws.Range("J2").Value = 71.833336
$ws.Range("K2").Value = 27150378
$ws.Range("L2").Value = 431.000016
$ws.Range("M2").Value = -27150265
$ws.Range("N2").Value = -657.000016
$ws.Range("H68").Value = 11912631
$ws.Range("I68").Value = 27778206
$ws.Range("J68").Value = 13450
$ws.Range("K68").Value = 83334618
$ws.Range("L68").Value = 40350
$ws.Range("M68").Value = -83333807
$ws.Range("N68").Value = -41972
$ws.Range("H71").Value = 11912631
$ws.Range("I71").Value = 27778206
$ws.Range("J71").Value = 13450
$ws.Range("K71").Value = 250003854
$ws.Range("L71").Value = 121050
$ws.Range("M71").Value = -249999798
$ws.Range("N71").Value = -129162
$ws.Range("H92").Value = 75725.5
$ws.Range("I92").Value = 300002
$ws.Range("J92").Value = 966.6667
$ws.Range("K92").Value = 900006
$ws.Range("L92").Value = 2900.0001
$ws.Range("M92").Value = -898758
$ws.Range("N92").Value = -5396.0001
$ws.Range("H107").Value = 487442.7
$ws.Range("I107").Value = 1648.25
$ws.Range("J107").Value = 973237.1
$ws.Range("K107").Value = 4944.75
$ws.Range("L107").Value = 2919711.3
$ws.Range("M107").Value = -3024.75
$ws.Range("N107").Value = -2923551.3
$ws.Range("H134").Value = 5008.2915
$ws.Range("I134").Value = 3619.9333
$ws.Range("J134").Value = 7322.222
$ws.Range("K134").Value = 10859.7999
$ws.Range("L134").Value = 21966.666
$ws.Range("M134").Value = -5789.7999
$ws.Range("N134").Value = -32106.666
$ws.Range("H137").Value = 4420.2
$ws.Range("I137").Value = 3490
$ws.Range("J137").Value = 6978.25
$ws.Range("K137").Value = 10470
$ws.Range("L137").Value = 20934.75
$ws.Range("M137").Value = -5370
$ws.Range("N137").Value = -31134.75
$ws.Range("H139").Value = 44827.824
$ws.Range("I139").Value = 48692.383
$ws.Range("J139").Value = 4250
$ws.Range("K139").Value = 146077.149
$ws.Range("L139").Value = 12750
$ws.Range("M139").Value = -140937.149
$ws.Range("N139").Value = -23030
$ws.Range("H140").Value = 253308.33
$ws.Range("I140").Value = 276154.53
$ws.Range("J140").Value = 2000
$ws.Range("K140").Value = 828463.5900000001
$ws.Range("L140").Value = 6000
$ws.Range("M140").Value = -823283.5900000001
$ws.Range("N140").Value = -16360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4441.88
$ws.Range("I70").Value = 4329.643
$ws.Range("J70").Value = 4584.727
$ws.Range("K70").Value = 4329.643
$ws.Range("L70").Value = 4584.727
$ws.Range("M70").Value = -4059.643
$ws.Range("N70").Value = -5124.727
$ws.Range("H73").Value = 4441.88
$ws.Range("I73").Value = 4329.643
$ws.Range("J73").Value = 4584.727
$ws.Range("K73").Value = 4329.643
$ws.Range("L73").Value = 4584.727
$ws.Range("M73").Value = -3393.643
$ws.Range("N73").Value = -6456.727
$ws.Range("H132").Value = 2219.12
$ws.Range("I132").Value = 1988.8182
$ws.Range("J132").Value = 2666.1765
$ws.Range("K132").Value = 5966.4546
$ws.Range("L132").Value = 7998.529500000001
$ws.Range("M132").Value = -3436.4546
$ws.Range("N132").Value = -13058.5295
$ws.Range("H136").Value = 11327.2
$ws.Range("J136").Value = 11327.2
$ws.Range("L136").Value = 33981.60000000001
$ws.Range("N136").Value = -39081.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2022163
$ws.Range("I40").Value = 5051255
$ws.Range("J40").Value = 2768.3333
$ws.Range("K40").Value = 5051255
$ws.Range("L40").Value = 2768.3333
$ws.Range("M40").Value = -5051119
$ws.Range("N40").Value = -3040.3333
$ws.Range("H122").Value = 100000
$ws.Range("I122").Value = 100000
$ws.Range("K122").Value = 300000
$ws.Range("M122").Value = -297550
$ws.Range("H132").Value = 1883.449
$ws.Range("I132").Value = 1690.9395
$ws.Range("J132").Value = 2280.5
$ws.Range("K132").Value = 5072.818499999999
$ws.Range("L132").Value = 6841.5
$ws.Range("M132").Value = -2542.818499999999
$ws.Range("N132").Value = -11901.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 781.42426
$ws.Range("I126").Value = 707
$ws.Range("J126").Value = 1116.3334
$ws.Range("K126").Value = 2121
$ws.Range("L126").Value = 3349.0002
$ws.Range("M126").Value = 349
$ws.Range("N126").Value = -8289.0002
